# Update "想去人数" (interested-count) figures in the 展览 and 全部类型 sheets.
# These values were refreshed by the data-generation job (gh-pages output).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 215
$wsExpo.Range("F6").Value  = 9831
$wsExpo.Range("F9").Value  = 1237
$wsExpo.Range("F10").Value = 3907
$wsExpo.Range("F12").Value = 114
$wsExpo.Range("F13").Value = 42
$wsExpo.Range("F16").Value = 543
$wsExpo.Range("F19").Value = 1444

# Sheet "全部类型" (all types) - same events, offset by one row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 215
$wsAll.Range("F7").Value  = 9831
$wsAll.Range("F10").Value = 1237
$wsAll.Range("F11").Value = 3907
$wsAll.Range("F13").Value = 114
$wsAll.Range("F14").Value = 42
$wsAll.Range("F17").Value = 543
$wsAll.Range("F20").Value = 1444
